$wb = $excel.ActiveWorkbook

# --- Restricciones_del_lider ---
$ws1 = $wb.Worksheets.Item("Restricciones_del_lider")

$ws1.Range("A2").Formula = '="1.0499999999999998 - x"'
$ws1.Range("A2").Copy()
$ws1.Range("A2").PasteSpecial(-4163)

$ws1.Range("B2").Formula = '="-2.05"'
$ws1.Range("B2").Copy()
$ws1.Range("B2").PasteSpecial(-4163)

$ws1.Range("D2").Formula = '="0.24"'
$ws1.Range("D2").Copy()
$ws1.Range("D2").PasteSpecial(-4163)

$ws1.Range("A3").Formula = '="-1.05 + x"'
$ws1.Range("A3").Copy()
$ws1.Range("A3").PasteSpecial(-4163)

$ws1.Range("B3").Formula = '="0.050000000000000044"'
$ws1.Range("B3").Copy()
$ws1.Range("B3").PasteSpecial(-4163)

$ws1.Range("D3").Formula = '="0.72"'
$ws1.Range("D3").Copy()
$ws1.Range("D3").PasteSpecial(-4163)

# --- Restricciones_del_follower ---
$ws2 = $wb.Worksheets.Item("Restricciones_del_follower")

$ws2.Range("A2").Formula = '="-2.85 + y"'
$ws2.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4163)

$ws2.Range("B2").Formula = '="1.85"'
$ws2.Range("B2").Copy()
$ws2.Range("B2").PasteSpecial(-4163)

$ws2.Range("D2").Formula = '="0.47"'
$ws2.Range("D2").Copy()
$ws2.Range("D2").PasteSpecial(-4163)

$ws2.Range("E2").Formula = '="0"'
$ws2.Range("E2").Copy()
$ws2.Range("E2").PasteSpecial(-4163)

$ws2.Range("F2").Formula = '="0"'
$ws2.Range("F2").Copy()
$ws2.Range("F2").PasteSpecial(-4163)

$ws2.Range("A3").Formula = '="2.85 - y"'
$ws2.Range("A3").Copy()
$ws2.Range("A3").PasteSpecial(-4163)

$ws2.Range("B3").Formula = '="-3.85"'
$ws2.Range("B3").Copy()
$ws2.Range("B3").PasteSpecial(-4163)

$ws2.Range("D3").Formula = '="0.88"'
$ws2.Range("D3").Copy()
$ws2.Range("D3").PasteSpecial(-4163)

$ws2.Range("F3").Formula = '="0"'
$ws2.Range("F3").Copy()
$ws2.Range("F3").PasteSpecial(-4163)

# --- Punto_modificado ---
$ws3 = $wb.Worksheets.Item("Punto_modificado")

$ws3.Range("A2").Formula = '="1.05"'
$ws3.Range("A2").Copy()
$ws3.Range("A2").PasteSpecial(-4163)

$ws3.Range("B2").Formula = '="2.85"'
$ws3.Range("B2").Copy()
$ws3.Range("B2").PasteSpecial(-4163)

# --- Vector_bf ---
$ws4 = $wb.Worksheets.Item("Vector_bf")

$ws4.Range("A2").Formula = '="-1.4248750000000001"'
$ws4.Range("A2").Copy()
$ws4.Range("A2").PasteSpecial(-4163)

# --- Vector_BF ---
$ws5 = $wb.Worksheets.Item("Vector_BF")

$ws5.Range("A2").Formula = '="-1.48"'
$ws5.Range("A2").Copy()
$ws5.Range("A2").PasteSpecial(-4163)

$ws5.Range("A3").Formula = '="1.0"'
$ws5.Range("A3").Copy()
$ws5.Range("A3").PasteSpecial(-4163)

$excel.CutCopyMode = 0

Write-Output "edits applied"
